# Weekly update: insert a new price record for "Locoto" (row shifts all
# subsequent records down by one) and backfill its columns.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 44, pushing the former row 44..87 down to 45..88.
$ws.Rows.Item(44).Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Range("A44").Value = 1
$ws.Range("B44").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C44").Value = "Arica y Parinacota"
$ws.Range("D44").Value = 44554
$ws.Range("E44").Value = 15
$ws.Range("F44").Value = 100112042
$ws.Range("G44").Value = "Locoto"
$ws.Range("H44").Value = "Sin especificar"
$ws.Range("I44").Value = "Primera"
$ws.Range("J44").Value = 160
$ws.Range("K44").Value = 10000
$ws.Range("L44").Value = 11000
$ws.Range("M44").Value = 10500
$ws.Range("N44").Value = "$/caja 20 kilos"
$ws.Range("O44").Value = "Región de Arica y Parinacota"
$ws.Range("P44").Value = 525
$ws.Range("Q44").Value = 20
$ws.Range("R44").Value = "Hortaliza"
